# #47 update myevent CD document again
#
# Slide 6 ("My Event Component" diagram):
#   - "events"        -> "myevents"        (+ widen the textbox)
#   - "<event-form>"  -> "<my-event-form>" (+ widen the textbox)
#   - "<event-card>"  -> "<my-event-card>" (+ widen the textbox)
#   - the red annotation mentioning "<event-form>" is also updated to
#     "<my-event-form>" (box size unchanged there)

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(6)

# --- "events[ ]" textbox -------------------------------------------------
$eventsBox = $s.Shapes.Item(5)
$eventsBox.TextFrame.TextRange.Runs(1).Text = "myevents"
$eventsBox.Width = 102.95074899999938

# --- "<event-form>" textbox ----------------------------------------------
$eventFormBox = $s.Shapes.Item(10)
$eventFormBox.TextFrame.TextRange.Runs(1).Text = "<my-event-form>"
$eventFormBox.Width = 160.69413699999916

# --- "<event-card>" textbox -----------------------------------------------
$eventCardBox = $s.Shapes.Item(11)
$eventCardBox.TextFrame.TextRange.Runs(1).Text = "<my-event-card>"
$eventCardBox.Width = 156.22090099999897

# --- red callout note referencing "<event-form>" --------------------------
$noteBox = $s.Shapes.Item(28)
$noteBox.TextFrame.TextRange.Paragraphs(2).Runs(1).Text = "- <my-event-form> is shown when we click add+"
